$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 4")

# --- Insert one new row at position 28. This pushes the existing "Today's
# Total" header section (old rows 29/30) and the two other week-summary
# sections (old rows 33/34 and 37/38) all down by one row each, leaving
# rows 27 and 28 free to be filled with the new Contact-Quality entry and a
# new "Today's Total" row for this block.
$ws.Rows.Item(28).Insert()

# --- Row 25: total time for that entry changed from 0.5 to 1 hour.
$ws.Range("D25").Value = 1

# --- Row 26: add the Finish time (B26) and the Total Time (D26) which were
# previously missing (the entry was left open/unfinished).
$ws.Range("B26").Value = 0.069444444444444434
$ws.Range("B26").NumberFormat = "h:mm"
$ws.Range("D26").Value = 0.66

# --- New row 27: a new Contact Quality entry.
$ws.Range("A27").Value = 0.069444444444444434
$ws.Range("A27").NumberFormat = "h:mm"
$ws.Range("B27").Value = 0.16666666666666666
$ws.Range("B27").NumberFormat = "h:mm"
$ws.Range("C27").Value = "Updating Contact Quality Images"
$ws.Range("D27").Value = 2.33

# --- New row 28: "Today's Total" summary row for this block, styled like
# the other "Today's Total" rows (centered, merged A:C).
$ws.Range("A28").Value = "Today's Total"
$ws.Range("A28:C28").HorizontalAlignment = -4108
$ws.Range("A28:C28").Merge()
$ws.Range("D28").Value = 9

# --- New row 42: "Week 4 Total" summary row at the very bottom, summing the
# four "Today's Total" cells for the week.
$ws.Range("A42").Value = "Week 4 Total"
$ws.Range("A42:C42").HorizontalAlignment = -4108
$ws.Range("A42:C42").Merge()
$ws.Range("D42").Formula = "=SUM(D28,D18,D12,D5)"

# --- Restore the view to what the saved workbook shows: scrolled near the
# top, with the newly-added blank row below the Week 4 Total selected.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("A43").Select()
